$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume 1h) updates per row.
# Numeric-looking Price strings are prefixed with an apostrophe so Excel
# stores them as text (matching the original inlineStr cell type) instead
# of auto-converting them to numbers.

$ws.Range("D2").Value = "27.048.37"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.621.45"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'213.79"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").Value = "'19.90"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "1.604.65"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "'0.538"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "27.036.90"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "'64.41"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "'214.87"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'6.81"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'4.34"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").Value = "'2.34"
$ws.Range("E23").Value = "  -6.51%  "
$ws.Range("D24").Value = "'9.01"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "'147.38"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").Value = "'0.0510"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").Value = "'0.722"
$ws.Range("E33").Value = "  +32.67%  "
$ws.Range("D34").Value = "'2.99"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").Value = "1.334.84"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").Value = "'1.56"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").Value = "'0.839"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'2.22"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'0.794"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "'63.85"
$ws.Range("E44").Value = "  +3.56%  "
$ws.Range("D45").Value = "1.760.16"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D48").Value = "'0.854"
$ws.Range("E48").Value = "  +27.78%  "
$ws.Range("D49").Value = "'0.0996"
$ws.Range("E49").Value = "  +4.08%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "'7.57"
$ws.Range("E51").Value = "  -1.09%  "

Write-Output "Updated cryptos list"
